# Apply scheduled runner updates to leve-profit sheets.
# Values were recomputed upstream (market price refresh); this script
# writes the refreshed currentAveragePrice / Leve price / profit figures
# back into each affected worksheet cell, matching the source diff.

$wb = $excel.ActiveWorkbook


# --- ALC sheet ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 174.25
$ws.Cells.Item(6, 9).Value = 199
$ws.Cells.Item(6, 11).Value = 597
$ws.Cells.Item(6, 13).Value = -485
$ws.Cells.Item(8, 8).Value = 647.2
$ws.Cells.Item(8, 9).Value = 647.2
$ws.Cells.Item(8, 11).Value = 1941.6
$ws.Cells.Item(8, 13).Value = -1802.6
$ws.Cells.Item(11, 8).Value = 302.75
$ws.Cells.Item(11, 9).Value = 302.75
$ws.Cells.Item(11, 11).Value = 302.75
$ws.Cells.Item(11, 13).Value = -162.75
$ws.Cells.Item(17, 8).Value = 3564.8572
$ws.Cells.Item(17, 10).Value = 3564.8572
$ws.Cells.Item(17, 12).Value = 10694.5716
$ws.Cells.Item(17, 14).Value = -11030.5716
$ws.Cells.Item(40, 8).Value = 2177.4167
$ws.Cells.Item(40, 9).Value = 2100
$ws.Cells.Item(40, 11).Value = 2100
$ws.Cells.Item(40, 13).Value = -1925
$ws.Cells.Item(43, 8).Value = 4332
$ws.Cells.Item(43, 9).Value = 6664
$ws.Cells.Item(43, 11).Value = 6664
$ws.Cells.Item(43, 13).Value = -6595
$ws.Cells.Item(92, 8).Value = 365.13333
$ws.Cells.Item(92, 9).Value = 347.41666
$ws.Cells.Item(92, 10).Value = 436
$ws.Cells.Item(92, 11).Value = 347.41666
$ws.Cells.Item(92, 12).Value = 436
$ws.Cells.Item(92, 13).Value = 900.58334
$ws.Cells.Item(92, 14).Value = -2932
$ws.Cells.Item(107, 8).Value = 282
$ws.Cells.Item(107, 9).Value = 290.3
$ws.Cells.Item(107, 11).Value = 290.3
$ws.Cells.Item(107, 13).Value = 1629.7
$ws.Cells.Item(112, 8).Value = 1193.0571
$ws.Cells.Item(112, 10).Value = 1588.9546
$ws.Cells.Item(112, 12).Value = 4766.8638
$ws.Cells.Item(112, 14).Value = -6982.8638
$ws.Cells.Item(116, 8).Value = 6981.5
$ws.Cells.Item(116, 10).Value = 6971.25
$ws.Cells.Item(116, 12).Value = 6971.25
$ws.Cells.Item(116, 14).Value = -13855.25
$ws.Cells.Item(125, 8).Value = 2014.5
$ws.Cells.Item(125, 9).Value = 1999
$ws.Cells.Item(125, 11).Value = 17991
$ws.Cells.Item(125, 13).Value = -15531
$ws.Cells.Item(138, 8).Value = 2733.9805
$ws.Cells.Item(138, 10).Value = 3167.9167
$ws.Cells.Item(138, 12).Value = 9503.750100000001
$ws.Cells.Item(138, 14).Value = -19783.7501

# --- ARM sheet ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 6483.8237
$ws.Cells.Item(32, 9).Value = 6483.8237
$ws.Cells.Item(32, 11).Value = 6483.8237
$ws.Cells.Item(32, 13).Value = -6196.8237
$ws.Cells.Item(45, 8).Value = 1937.5
$ws.Cells.Item(45, 9).Value = 1937.5
$ws.Cells.Item(45, 11).Value = 1937.5
$ws.Cells.Item(45, 13).Value = -1560.5
$ws.Cells.Item(132, 8).Value = 3152.8823
$ws.Cells.Item(132, 9).Value = 2450.375
$ws.Cells.Item(132, 11).Value = 7351.125
$ws.Cells.Item(132, 13).Value = -4821.125

# --- BSM sheet ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(107, 8).Value = 2140.3333
$ws.Cells.Item(107, 9).Value = 2140.3333
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 2140.3333
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -220.3332999999998
$ws.Cells.Item(107, 14).ClearContents()

# --- CRP sheet ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = 1146.8334
$ws.Cells.Item(16, 9).Value = 1247.625
$ws.Cells.Item(16, 10).Value = 945.25
$ws.Cells.Item(16, 11).Value = 1247.625
$ws.Cells.Item(16, 12).Value = 945.25
$ws.Cells.Item(16, 13).Value = -960.625
$ws.Cells.Item(16, 14).Value = -1519.25
$ws.Cells.Item(31, 8).Value = 5868.6665
$ws.Cells.Item(31, 9).Value = 1566.6666
$ws.Cells.Item(31, 11).Value = 1566.6666
$ws.Cells.Item(31, 13).Value = -1271.6666
$ws.Cells.Item(34, 8).Value = 5868.6665
$ws.Cells.Item(34, 9).Value = 1566.6666
$ws.Cells.Item(34, 11).Value = 1566.6666
$ws.Cells.Item(34, 13).Value = -1364.6666
$ws.Cells.Item(58, 8).Value = 2719.3333
$ws.Cells.Item(58, 9).Value = 2262
$ws.Cells.Item(58, 11).Value = 2262
$ws.Cells.Item(58, 13).Value = -2059
$ws.Cells.Item(107, 8).Value = 1682.7693
$ws.Cells.Item(107, 9).Value = 782.2857
$ws.Cells.Item(107, 11).Value = 782.2857
$ws.Cells.Item(107, 13).Value = 1137.7143
$ws.Cells.Item(113, 8).Value = 1146.8334
$ws.Cells.Item(113, 9).Value = 1247.625
$ws.Cells.Item(113, 10).Value = 945.25
$ws.Cells.Item(113, 11).Value = 1247.625
$ws.Cells.Item(113, 12).Value = 945.25
$ws.Cells.Item(113, 13).Value = 922.375
$ws.Cells.Item(113, 14).Value = -5285.25
$ws.Cells.Item(122, 8).Value = 3121
$ws.Cells.Item(122, 9).Value = 2532.6667
$ws.Cells.Item(122, 10).Value = 3373.1428
$ws.Cells.Item(122, 11).Value = 7598.000100000001
$ws.Cells.Item(122, 12).Value = 10119.4284
$ws.Cells.Item(122, 13).Value = -5148.000100000001
$ws.Cells.Item(122, 14).Value = -15019.4284
$ws.Cells.Item(136, 8).Value = 2719.3333
$ws.Cells.Item(136, 9).Value = 2262
$ws.Cells.Item(136, 11).Value = 6786
$ws.Cells.Item(136, 13).Value = -4236

# --- CUL sheet ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(14, 8).Value = 430.66666
$ws.Cells.Item(14, 9).Value = 430.66666
$ws.Cells.Item(14, 11).Value = 1291.99998
$ws.Cells.Item(14, 13).Value = -1118.99998
$ws.Cells.Item(56, 8).Value = 17635.639
$ws.Cells.Item(56, 9).Value = 17635.639
$ws.Cells.Item(56, 11).Value = 17635.639
$ws.Cells.Item(56, 13).Value = -17105.639
$ws.Cells.Item(64, 8).Value = 950
$ws.Cells.Item(64, 9).Value = 1000
$ws.Cells.Item(64, 10).Value = 900
$ws.Cells.Item(64, 11).Value = 3000
$ws.Cells.Item(64, 12).Value = 2700
$ws.Cells.Item(64, 13).Value = -2730
$ws.Cells.Item(64, 14).Value = -3240
$ws.Cells.Item(67, 8).Value = 950
$ws.Cells.Item(67, 9).Value = 1000
$ws.Cells.Item(67, 10).Value = 900
$ws.Cells.Item(67, 11).Value = 3000
$ws.Cells.Item(67, 12).Value = 2700
$ws.Cells.Item(67, 13).Value = -2064
$ws.Cells.Item(67, 14).Value = -4572
$ws.Cells.Item(75, 8).Value = 995.2
$ws.Cells.Item(75, 10).Value = 1325.3334
$ws.Cells.Item(75, 12).Value = 3976.0002
$ws.Cells.Item(75, 14).Value = -5972.0002
$ws.Cells.Item(78, 8).Value = 995.2
$ws.Cells.Item(78, 10).Value = 1325.3334
$ws.Cells.Item(78, 12).Value = 11928.0006
$ws.Cells.Item(78, 14).Value = -21912.0006
$ws.Cells.Item(113, 8).Value = 1768.7059
$ws.Cells.Item(113, 10).Value = 1687.9333
$ws.Cells.Item(113, 12).Value = 5063.7999
$ws.Cells.Item(113, 14).Value = -9403.7999

# --- LTW sheet ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 4872.1665
$ws.Cells.Item(40, 9).Value = 4872.1665
$ws.Cells.Item(40, 11).Value = 4872.1665
$ws.Cells.Item(40, 13).Value = -4736.1665
$ws.Cells.Item(61, 8).Value = 1758.4
$ws.Cells.Item(61, 9).Value = 948.25
$ws.Cells.Item(61, 11).Value = 948.25
$ws.Cells.Item(61, 13).Value = -746.25
$ws.Cells.Item(113, 8).Value = 1758.4
$ws.Cells.Item(113, 9).Value = 948.25
$ws.Cells.Item(113, 11).Value = 948.25
$ws.Cells.Item(113, 13).Value = 1221.75

# --- WVR sheet ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(113, 8).Value = 1008.1539
$ws.Cells.Item(113, 9).Value = 1046.5454
$ws.Cells.Item(113, 10).Value = 797
$ws.Cells.Item(113, 11).Value = 3139.6362
$ws.Cells.Item(113, 12).Value = 2391
$ws.Cells.Item(113, 13).Value = -969.6361999999999
$ws.Cells.Item(113, 14).Value = -6731
$ws.Cells.Item(132, 8).Value = 3056.0417
$ws.Cells.Item(132, 9).Value = 2746.9285
$ws.Cells.Item(132, 10).Value = 3488.8
$ws.Cells.Item(132, 11).Value = 8240.7855
$ws.Cells.Item(132, 12).Value = 10466.4
$ws.Cells.Item(132, 13).Value = -5710.7855
$ws.Cells.Item(132, 14).Value = -15526.4
$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()
